$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new date
$ws.Name = "Apr 14 2022"

# Update the date value in A2 (stored as Excel serial date 44665 = 2022-04-14)
$ws.Range("A2").Value = 44665

# Update the active selection on the sheet
$ws.Range("E9").Select()
